$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows before row 450, pushing existing rows 450:555 down to 453:558
$ws.Rows.Item(450).Resize(3).Insert()

# New row 450
$ws.Cells.Item(450, 1).Value = 3
$ws.Cells.Item(450, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(450, 3).Value = 'Coquimbo'
$ws.Cells.Item(450, 4).Value = 44943
$ws.Cells.Item(450, 5).Value = 5
$ws.Cells.Item(450, 6).Value = 100112027
$ws.Cells.Item(450, 7).Value = 'Melón'
$ws.Cells.Item(450, 8).Value = 'Tuna'
$ws.Cells.Item(450, 9).Value = 'Extra'
$ws.Cells.Item(450, 10).Value = 750
$ws.Cells.Item(450, 11).Value = 2000
$ws.Cells.Item(450, 12).Value = 2000
$ws.Cells.Item(450, 13).Value = 2000
$ws.Cells.Item(450, 14).Value = '$/unidad'
$ws.Cells.Item(450, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(450, 16).Value = 2000
$ws.Cells.Item(450, 17).Value = 1
$ws.Cells.Item(450, 18).Value = 'Hortaliza'

# New row 451
$ws.Cells.Item(451, 1).Value = 3
$ws.Cells.Item(451, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(451, 3).Value = 'Coquimbo'
$ws.Cells.Item(451, 4).Value = 44943
$ws.Cells.Item(451, 5).Value = 5
$ws.Cells.Item(451, 6).Value = 100112027
$ws.Cells.Item(451, 7).Value = 'Melón'
$ws.Cells.Item(451, 8).Value = 'Tuna'
$ws.Cells.Item(451, 9).Value = 'Primera'
$ws.Cells.Item(451, 10).Value = 760
$ws.Cells.Item(451, 11).Value = 1500
$ws.Cells.Item(451, 12).Value = 1500
$ws.Cells.Item(451, 13).Value = 1500
$ws.Cells.Item(451, 14).Value = '$/unidad'
$ws.Cells.Item(451, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(451, 16).Value = 1500
$ws.Cells.Item(451, 17).Value = 1
$ws.Cells.Item(451, 18).Value = 'Hortaliza'

# New row 452
$ws.Cells.Item(452, 1).Value = 3
$ws.Cells.Item(452, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(452, 3).Value = 'Coquimbo'
$ws.Cells.Item(452, 4).Value = 44943
$ws.Cells.Item(452, 5).Value = 5
$ws.Cells.Item(452, 6).Value = 100112027
$ws.Cells.Item(452, 7).Value = 'Melón'
$ws.Cells.Item(452, 8).Value = 'Tuna'
$ws.Cells.Item(452, 9).Value = 'Segunda'
$ws.Cells.Item(452, 10).Value = 750
$ws.Cells.Item(452, 11).Value = 1000
$ws.Cells.Item(452, 12).Value = 1000
$ws.Cells.Item(452, 13).Value = 1000
$ws.Cells.Item(452, 14).Value = '$/unidad'
$ws.Cells.Item(452, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(452, 16).Value = 1000
$ws.Cells.Item(452, 17).Value = 1
$ws.Cells.Item(452, 18).Value = 'Hortaliza'

# Apply the date number format to column D of the new rows, matching the existing date cells
$ws.Range('D450:D452').NumberFormat = $ws.Range('D449').NumberFormat
